$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 only held the now-redundant label "grandes regiões e unidades da
# federação" (no data in B:H). Delete the entire row so every region row
# below it shifts up by one, aligning the region totals correctly and
# dropping the trailing empty row (old row 38).
$ws.Rows(6).Delete()
